# Update column F ("dSF") values on Sheet1 to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -1
    3  = 3
    4  = -2
    5  = -1
    7  = 2
    8  = 7
    9  = -4
    10 = -1
    12 = 1
    13 = -1
    15 = -5
    16 = 2
    17 = 1
    18 = -1
    19 = 1
    21 = 6
    22 = 0
    23 = -1
    24 = -3
    25 = -10
    27 = 10
    28 = -5
    29 = -1
    31 = -1
    34 = 2
    35 = -3
    37 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

$wb.Save()
